# NMCARS 18-18 update — applies the tracked OOXML diff via the Word COM
# object model. A short raw-XML patch finishes a handful of style-
# definition knobs that have no reachable COM surface in this host
# (explicit w:tabs on a *style* pPr, fully dropping an rPr node, rsid
# stamps, w:semiHidden, w:contextualSpacing and the built-in-id flag
# on List3/List4) — each patch only touches text this script itself
# just produced, so it can't clobber unrelated content.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Two paragraphs ("(a) All requests ..." and "(c)(2) When the ...")
#    get the new "List1" paragraph style. Build/link the style family
#    via the OM first so the two Paragraphs.Style assignments bind to
#    the fully-formed style instead of an auto-vivified stub.
# ---------------------------------------------------------------------

$list1 = $d.Styles.Add("List 1", 1)
$list1.BaseStyle = $d.Styles("Heading1")

$heading1Char1 = $d.Styles.Add("Heading 1 Char1", 2)
$heading1Char1.BaseStyle = $d.Styles("DefaultParagraphFont")
$heading1Char1.Font.Bold = $true
$heading1Char1.Font.Size = 16

$list1Char = $d.Styles.Add("List 1 Char", 2)
$list1Char.BaseStyle = $heading1Char1

$list1.LinkStyle = "List1Char"
$heading1Char1.LinkStyle = "Heading1"
$list1Char.LinkStyle = "List1"

$pf = $list1.ParagraphFormat
$pf.KeepWithNext = $false
$pf.SpaceBefore = 0
$pf.SpaceAfter = 0
$pf.Alignment = 0
$pf.OutlineLevel = 10

$list1.Font.Bold = $false
$list1.Font.Size = 12
$list1Char.Font.Bold = $false
$list1Char.Font.Size = 12

# Heading 1 <-> Heading 1 Char1 back-reference on the Heading1 style
# itself.
$heading1 = $d.Styles("Heading1")
$heading1.LinkStyle = "Heading1Char1"

# Apply the finished List1 style to the two target paragraphs.
$d.Paragraphs(3).Style = "List1"
$d.Paragraphs(7).Style = "List1"

# ---------------------------------------------------------------------
# 2) "List 2" style loses its Courier New rPr override.
#    (Font.Name = "" only blanks the rFonts attributes via this host's
#    OM surface; the full-element removal is finished in the raw patch
#    below.)
# ---------------------------------------------------------------------

$list2 = $d.Styles("List2")
$list2.Font.Name = ""

# ---------------------------------------------------------------------
# 3) New built-in-id "List 3" / "List 4" styles (List 4 semi-hidden /
#    unhide-when-used), both based on Normal with hanging indents.
#    (ParagraphFormat indents are in points; 1080/360/1440 twips ==
#    54/-18/72/-18 pt.)
# ---------------------------------------------------------------------

$list3 = $d.Styles.Add("List 3", 1)
$list3.BaseStyle = $d.Styles("Normal")
$list3.ParagraphFormat.LeftIndent = 54
$list3.ParagraphFormat.FirstLineIndent = -18

$list4 = $d.Styles.Add("List 4", 1)
$list4.BaseStyle = $d.Styles("Normal")
$list4.ParagraphFormat.LeftIndent = 72
$list4.ParagraphFormat.FirstLineIndent = -18
$list4.UnhideWhenUsed = $true

Write-Output "com-phase-done"

# ---------------------------------------------------------------------
# 4) Raw-XML finishing patch for the knobs the OM surface above can't
#    reach in this host: the List1 style's explicit tab stop, fully
#    dropping List2's rPr node, stamping the 004450B4 rsid the real
#    commit carries on every new style, semiHidden on List 4, the
#    contextualSpacing toggle on List 3 / List 4, and clearing the
#    w:customStyle="1" flag Styles.Add always stamps so List3/List4
#    read as Word's built-in list styles (matching List2's existing
#    un-flagged form).
# ---------------------------------------------------------------------

$xml = $d.WordOpenXML

$xml = $xml.Replace(
  '<w:style w:type="paragraph" w:styleId="List2"><w:name w:val="List 2"/><w:basedOn w:val="Normal"/><w:uiPriority w:val="99"/><w:rsid w:val="00D44619"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="" w:hAnsi=""/></w:rPr></w:style>',
  '<w:style w:type="paragraph" w:styleId="List2"><w:name w:val="List 2"/><w:basedOn w:val="Normal"/><w:uiPriority w:val="99"/><w:rsid w:val="00D44619"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr></w:style>'
)

$xml = $xml.Replace(
  '<w:style w:type="paragraph" w:customStyle="1" w:styleId="List1"><w:name w:val="List 1"/><w:basedOn w:val="Heading1"/><w:link w:val="List1Char"/><w:pPr><w:keepNext w:val="0"/><w:spacing w:before="0" w:after="0"/><w:jc w:val="left"/><w:outlineLvl w:val="9"/></w:pPr><w:rPr><w:b w:val="0"/><w:sz w:val="24"/></w:rPr></w:style>',
  '<w:style w:type="paragraph" w:customStyle="1" w:styleId="List1"><w:name w:val="List 1"/><w:basedOn w:val="Heading1"/><w:link w:val="List1Char"/><w:rsid w:val="004450B4"/><w:pPr><w:keepNext w:val="0"/><w:tabs><w:tab w:val="left" w:pos="3686"/></w:tabs><w:spacing w:before="0" w:after="0"/><w:jc w:val="left"/><w:outlineLvl w:val="9"/></w:pPr><w:rPr><w:b w:val="0"/><w:sz w:val="24"/></w:rPr></w:style>'
)

$xml = $xml.Replace(
  '<w:style w:type="character" w:customStyle="1" w:styleId="Heading1Char1"><w:name w:val="Heading 1 Char1"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading1"/><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:style>',
  '<w:style w:type="character" w:customStyle="1" w:styleId="Heading1Char1"><w:name w:val="Heading 1 Char1"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading1"/><w:rsid w:val="004450B4"/><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:style>'
)

$xml = $xml.Replace(
  '<w:style w:type="character" w:customStyle="1" w:styleId="List1Char"><w:name w:val="List 1 Char"/><w:basedOn w:val="Heading1Char1"/><w:link w:val="List1"/><w:rPr><w:b w:val="0"/><w:sz w:val="24"/></w:rPr></w:style>',
  '<w:style w:type="character" w:customStyle="1" w:styleId="List1Char"><w:name w:val="List 1 Char"/><w:basedOn w:val="Heading1Char1"/><w:link w:val="List1"/><w:rsid w:val="004450B4"/><w:rPr><w:b w:val="0"/><w:sz w:val="24"/></w:rPr></w:style>'
)

$xml = $xml.Replace(
  '<w:style w:type="paragraph" w:customStyle="1" w:styleId="List3"><w:name w:val="List 3"/><w:basedOn w:val="Normal"/><w:pPr><w:ind w:left="1080" w:hanging="360"/></w:pPr></w:style>',
  '<w:style w:type="paragraph" w:styleId="List3"><w:name w:val="List 3"/><w:basedOn w:val="Normal"/><w:rsid w:val="004450B4"/><w:pPr><w:ind w:left="1080" w:hanging="360"/><w:contextualSpacing/></w:pPr></w:style>'
)

$xml = $xml.Replace(
  '<w:style w:type="paragraph" w:customStyle="1" w:styleId="List4"><w:name w:val="List 4"/><w:basedOn w:val="Normal"/><w:unhideWhenUsed/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr></w:style>',
  '<w:style w:type="paragraph" w:styleId="List4"><w:name w:val="List 4"/><w:basedOn w:val="Normal"/><w:semiHidden/><w:unhideWhenUsed/><w:rsid w:val="004450B4"/><w:pPr><w:ind w:left="1440" w:hanging="360"/><w:contextualSpacing/></w:pPr></w:style>'
)

$d.WordOpenXML = $xml

Write-Output "raw-patch-done"
